# Update "想去人数" (attendance count) figures on the "展览" (Exhibition)
# and "全部类型" (All types) worksheets.
#
# 展览  sheet: F2 159->160, F3 61->62, F4 264->265, F5 4036->4044, F7 446->447
# 全部类型 sheet: F2 159->160, F3 61->62, F4 264->265, F5 4036->4044, F9 446->447

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 160
$ws1.Range("F3").Value = 62
$ws1.Range("F4").Value = 265
$ws1.Range("F5").Value = 4044
$ws1.Range("F7").Value = 447

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 160
$ws4.Range("F3").Value = 62
$ws4.Range("F4").Value = 265
$ws4.Range("F5").Value = 4044
$ws4.Range("F9").Value = 447
